# Generate Report for Handoff
# Rename the source markdown file from
#   63648775-0410-4ac6-8bdc-55c61680613c.md
# to
#   60fbb8d4-661e-4a21-9688-2ae0bcc04d4b.md
# across all three report sheets, and bump the related handoff
# timestamps / generated xliff file names.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# File name / path updates (shared across sheets).
$wsOverview.Range("A2").Value = "60fbb8d4-661e-4a21-9688-2ae0bcc04d4b.md"
$wsOverview.Range("B2").Value = "e2e\60fbb8d4-661e-4a21-9688-2ae0bcc04d4b.md"

$wsZhCn.Range("A2").Value = "60fbb8d4-661e-4a21-9688-2ae0bcc04d4b.md"
$wsDeDe.Range("A2").Value = "60fbb8d4-661e-4a21-9688-2ae0bcc04d4b.md"

# Latest HO Xliff Generate Date (Overview!G2) is the same underlying value
# as the Latest Handback DateTime on the de-de sheet (de-de!H2) - bump both.
$wsOverview.Range("G2").Value = "2016-08-16 00:52:56"
$wsDeDe.Range("H2").Value = "2016-08-16 00:52:56"

# Latest Handoff File names (new content hash with the new source file name).
$wsZhCn.Range("G2").Value = "60fbb8d4-661e-4a21-9688-2ae0bcc04d4b.e9ffd5825af93b237f144eebeaf01e3b05938745.zh-cn.xlf"
$wsDeDe.Range("G2").Value = "60fbb8d4-661e-4a21-9688-2ae0bcc04d4b.e9ffd5825af93b237f144eebeaf01e3b05938745.de-de.xlf"

# Latest Handoff Datetime on zh-cn sheet.
$wsZhCn.Range("H2").Value = "2016-08-16 00:52:51"
